$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.026400566431506
$ws.Cells.Item(2, 4).Value = 1.034032452575097
$ws.Cells.Item(2, 5).Value = 1.026631291024233
$ws.Cells.Item(2, 6).Value = 1.041131495943915
$ws.Cells.Item(2, 9).Value = 1.029884742224066
$ws.Cells.Item(2, 10).Value = 1.031564217223217
$ws.Cells.Item(2, 11).Value = 1.036833092325719
$ws.Cells.Item(2, 12).Value = 1.029453373643026
$ws.Cells.Item(2, 13).Value = 1.043911879158821
$ws.Cells.Item(2, 14).Value = 1.014540865664894

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.027272379534459
$ws.Cells.Item(3, 4).Value = 1.034831428721074
$ws.Cells.Item(3, 5).Value = 1.027369414943965
$ws.Cells.Item(3, 6).Value = 1.042088092640608
$ws.Cells.Item(3, 9).Value = 1.029956334752582
$ws.Cells.Item(3, 10).Value = 1.032076174201258
$ws.Cells.Item(3, 11).Value = 1.037441199128077
$ws.Cells.Item(3, 12).Value = 1.029999240676111
$ws.Cells.Item(3, 13).Value = 1.044678657292176
$ws.Cells.Item(3, 14).Value = 1.014711986659328

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.027837009318187
$ws.Cells.Item(4, 4).Value = 1.035349252805494
$ws.Cells.Item(4, 5).Value = 1.027847856257405
$ws.Cells.Item(4, 6).Value = 1.042708261488899
$ws.Cells.Item(4, 9).Value = 1.030001226733613
$ws.Cells.Item(4, 10).Value = 1.032407321195355
$ws.Cells.Item(4, 11).Value = 1.037834864291048
$ws.Cells.Item(4, 12).Value = 1.030352613223806
$ws.Cells.Item(4, 13).Value = 1.045175379219029
$ws.Cells.Item(4, 14).Value = 1.014822625204805

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.028074499956901
$ws.Cells.Item(5, 4).Value = 1.035567143824146
$ws.Cells.Item(5, 5).Value = 1.028049188985686
$ws.Cells.Item(5, 6).Value = 1.04296926291471
$ws.Cells.Item(5, 9).Value = 1.030019755601217
$ws.Cells.Item(5, 10).Value = 1.032546504633057
$ws.Cells.Item(5, 11).Value = 1.038000402731597
$ws.Cells.Item(5, 12).Value = 1.030501208066205
$ws.Cells.Item(5, 13).Value = 1.045384334925596
$ws.Cells.Item(5, 14).Value = 1.014869116086939

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.028114382691588
$ws.Cells.Item(6, 4).Value = 1.035603740233915
$ws.Cells.Item(6, 5).Value = 1.028083005080948
$ws.Cells.Item(6, 6).Value = 1.043013102716152
$ws.Cells.Item(6, 9).Value = 1.030022846505706
$ws.Cells.Item(6, 10).Value = 1.032569872305646
$ws.Cells.Item(6, 11).Value = 1.038028199764943
$ws.Cells.Item(6, 12).Value = 1.030526159927373
$ws.Cells.Item(6, 13).Value = 1.045419427311136
$ws.Cells.Item(6, 14).Value = 1.014876920833047

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.027840182208912
$ws.Cells.Item(7, 4).Value = 1.035352163500932
$ws.Cells.Item(7, 5).Value = 1.027850545707158
$ws.Cells.Item(7, 6).Value = 1.042711747895033
$ws.Cells.Item(7, 9).Value = 1.030001475668608
$ws.Cells.Item(7, 10).Value = 1.032409181093649
$ws.Cells.Item(7, 11).Value = 1.037837076061994
$ws.Cells.Item(7, 12).Value = 1.03035459861152
$ws.Cells.Item(7, 13).Value = 1.045178170771515
$ws.Cells.Item(7, 14).Value = 1.014823246503541

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.026695093617157
$ws.Cells.Item(8, 4).Value = 1.034302296802897
$ws.Cells.Item(8, 5).Value = 1.02688057145773
$ws.Cells.Item(8, 6).Value = 1.04145453552355
$ws.Cells.Item(8, 9).Value = 1.029909233631144
$ws.Cells.Item(8, 10).Value = 1.031737260301813
$ws.Cells.Item(8, 11).Value = 1.037038566913707
$ws.Cells.Item(8, 12).Value = 1.029637818036755
$ws.Cells.Item(8, 13).Value = 1.044170897486751
$ws.Cells.Item(8, 14).Value = 1.014598714792946

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.024681258565475
$ws.Cells.Item(9, 4).Value = 1.032458751847278
$ws.Cells.Item(9, 5).Value = 1.025177752198166
$ws.Cells.Item(9, 6).Value = 1.039248330839986
$ws.Cells.Item(9, 9).Value = 1.029735742232992
$ws.Cells.Item(9, 10).Value = 1.030552358579559
$ws.Cells.Item(9, 11).Value = 1.035632925071636
$ws.Cells.Item(9, 12).Value = 1.028376042465975
$ws.Cells.Item(9, 13).Value = 1.042400348429325
$ws.Cells.Item(9, 14).Value = 1.014202405821839

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.023341447639096
$ws.Cells.Item(10, 4).Value = 1.031234160691762
$ws.Cells.Item(10, 5).Value = 1.024046942640372
$ws.Cells.Item(10, 6).Value = 1.037783789272259
$ws.Cells.Item(10, 9).Value = 1.029612756355588
$ws.Cells.Item(10, 10).Value = 1.029761896901536
$ws.Cells.Item(10, 11).Value = 1.034696873915861
$ws.Cells.Item(10, 12).Value = 1.027535798831327
$ws.Cells.Item(10, 13).Value = 1.041223026660876
$ws.Cells.Item(10, 14).Value = 1.013937786002679

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.022761963723839
$ws.Cells.Item(11, 4).Value = 1.030704972910832
$ws.Cells.Item(11, 5).Value = 1.023558355244561
$ws.Cells.Item(11, 6).Value = 1.037151132757426
$ws.Cells.Item(11, 9).Value = 1.02955777241367
$ws.Cells.Item(11, 10).Value = 1.029419508363509
$ws.Cells.Item(11, 11).Value = 1.03429181765108
$ws.Cells.Item(11, 12).Value = 1.027172203264515
$ws.Cells.Item(11, 13).Value = 1.040713974111341
$ws.Cells.Item(11, 14).Value = 1.013823110175228

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.022546818618036
$ws.Cells.Item(12, 4).Value = 1.030508571183049
$ws.Cells.Item(12, 5).Value = 1.023377033189978
$ws.Cells.Item(12, 6).Value = 1.036916362948433
$ws.Cells.Item(12, 9).Value = 1.029537089576438
$ws.Cells.Item(12, 10).Value = 1.029292314182936
$ws.Cells.Item(12, 11).Value = 1.034141402128782
$ws.Cells.Item(12, 12).Value = 1.027037184359854
$ws.Cells.Item(12, 13).Value = 1.040525001333464
$ws.Cells.Item(12, 14).Value = 1.013780500839205

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.022592963379513
$ws.Cells.Item(13, 4).Value = 1.03055069265457
$ws.Cells.Item(13, 5).Value = 1.023415920070449
$ws.Cells.Item(13, 6).Value = 1.036966711570583
$ws.Cells.Item(13, 9).Value = 1.029541537841923
$ws.Cells.Item(13, 10).Value = 1.029319598463531
$ws.Cells.Item(13, 11).Value = 1.034173664918809
$ws.Cells.Item(13, 12).Value = 1.027066144692603
$ws.Cells.Item(13, 13).Value = 1.040565531543307
$ws.Cells.Item(13, 14).Value = 1.013789641295895

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.022744177693325
$ws.Cells.Item(14, 4).Value = 1.030688734964684
$ws.Cells.Item(14, 5).Value = 1.02354336381151
$ws.Cells.Item(14, 6).Value = 1.037131721971509
$ws.Cells.Item(14, 9).Value = 1.029556068051163
$ws.Cells.Item(14, 10).Value = 1.029408994768869
$ws.Cells.Item(14, 11).Value = 1.034279383425002
$ws.Cells.Item(14, 12).Value = 1.027161041803997
$ws.Cells.Item(14, 13).Value = 1.040698351269048
$ws.Cells.Item(14, 14).Value = 1.013819588349625

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.022837359251231
$ws.Cells.Item(15, 4).Value = 1.03077380892809
$ws.Cells.Item(15, 5).Value = 1.023621907499633
$ws.Cells.Item(15, 6).Value = 1.037233420454301
$ws.Cells.Item(15, 9).Value = 1.029564986243692
$ws.Cells.Item(15, 10).Value = 1.029464072721145
$ws.Cells.Item(15, 11).Value = 1.034344525480508
$ws.Cells.Item(15, 12).Value = 1.027219515945412
$ws.Cells.Item(15, 13).Value = 1.040780200772182
$ws.Cells.Item(15, 14).Value = 1.013838037925433

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.023379920167769
$ws.Cells.Item(16, 4).Value = 1.031269303792638
$ws.Cells.Item(16, 5).Value = 1.024079391068842
$ws.Cells.Item(16, 6).Value = 1.037825808341151
$ws.Cells.Item(16, 9).Value = 1.029616369055383
$ws.Cells.Item(16, 10).Value = 1.029784617814647
$ws.Cells.Item(16, 11).Value = 1.03472376175957
$ws.Cells.Item(16, 12).Value = 1.027559934543811
$ws.Cells.Item(16, 13).Value = 1.041256826426932
$ws.Cells.Item(16, 14).Value = 1.013945394723911

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.023720432806475
$ws.Cells.Item(17, 4).Value = 1.031580401823241
$ws.Cells.Item(17, 5).Value = 1.024366643590362
$ws.Cells.Item(17, 6).Value = 1.03819780029888
$ws.Cells.Item(17, 9).Value = 1.02964813732107
$ws.Cells.Item(17, 10).Value = 1.029985657770503
$ws.Cells.Item(17, 11).Value = 1.034961717228879
$ws.Cells.Item(17, 12).Value = 1.027773534097992
$ws.Cells.Item(17, 13).Value = 1.041555999090831
$ws.Cells.Item(17, 14).Value = 1.014012712038719

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.023919111948618
$ws.Cells.Item(18, 4).Value = 1.031761963001058
$ws.Cells.Item(18, 5).Value = 1.024534295344343
$ws.Cells.Item(18, 6).Value = 1.038414921495481
$ws.Cells.Item(18, 9).Value = 1.029666500229352
$ws.Cells.Item(18, 10).Value = 1.030102909951584
$ws.Cells.Item(18, 11).Value = 1.035100537630768
$ws.Cells.Item(18, 12).Value = 1.02789814573653
$ws.Cells.Item(18, 13).Value = 1.041730572313065
$ws.Cells.Item(18, 14).Value = 1.014051968018759

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.023986867194393
$ws.Cells.Item(19, 4).Value = 1.031823888083389
$ws.Cells.Item(19, 5).Value = 1.024591477545713
$ws.Cells.Item(19, 6).Value = 1.038488978684883
$ws.Cells.Item(19, 9).Value = 1.029672733165063
$ws.Cells.Item(19, 10).Value = 1.03014288800872
$ws.Cells.Item(19, 11).Value = 1.035147876000528
$ws.Cells.Item(19, 12).Value = 1.027940638884277
$ws.Cells.Item(19, 13).Value = 1.041790109245549
$ws.Cells.Item(19, 14).Value = 1.014065351741245

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.023683892405252
$ws.Cells.Item(20, 4).Value = 1.031547013286993
$ws.Cells.Item(20, 5).Value = 1.024335813528054
$ws.Cells.Item(20, 6).Value = 1.038157874115424
$ws.Cells.Item(20, 9).Value = 1.029644746154403
$ws.Cells.Item(20, 10).Value = 1.029964089217392
$ws.Cells.Item(20, 11).Value = 1.034936184263215
$ws.Cells.Item(20, 12).Value = 1.027750614536941
$ws.Cells.Item(20, 13).Value = 1.041523893354208
$ws.Cells.Item(20, 14).Value = 1.014005490462909

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.022699646046828
$ws.Cells.Item(21, 4).Value = 1.03064808045462
$ws.Cells.Item(21, 5).Value = 1.023505830346454
$ws.Cells.Item(21, 6).Value = 1.037083124241651
$ws.Cells.Item(21, 9).Value = 1.029551796421319
$ws.Cells.Item(21, 10).Value = 1.029382670212461
$ws.Cells.Item(21, 11).Value = 1.034248250833645
$ws.Cells.Item(21, 12).Value = 1.027133095947006
$ws.Cells.Item(21, 13).Value = 1.040659236066358
$ws.Cells.Item(21, 14).Value = 1.013810770061297

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.02208139647598
$ws.Cells.Item(22, 4).Value = 1.030083824837104
$ws.Cells.Item(22, 5).Value = 1.02298491966089
$ws.Cells.Item(22, 6).Value = 1.036408700466659
$ws.Cells.Item(22, 9).Value = 1.029491854993712
$ws.Cells.Item(22, 10).Value = 1.029017018001763
$ws.Cells.Item(22, 11).Value = 1.033815954242299
$ws.Cells.Item(22, 12).Value = 1.026745050350977
$ws.Cells.Item(22, 13).Value = 1.040116240477362
$ws.Cells.Item(22, 14).Value = 1.013688262992526

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.022409086332849
$ws.Cells.Item(23, 4).Value = 1.030382857836031
$ws.Cells.Item(23, 5).Value = 1.023260975291673
$ws.Cells.Item(23, 6).Value = 1.036766100143251
$ws.Cells.Item(23, 9).Value = 1.02952377308486
$ws.Cells.Item(23, 10).Value = 1.029210865380427
$ws.Cells.Item(23, 11).Value = 1.034045100173764
$ws.Cells.Item(23, 12).Value = 1.026950740053507
$ws.Cells.Item(23, 13).Value = 1.040404030717127
$ws.Cells.Item(23, 14).Value = 1.013753213610234

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.023700403242463
$ws.Cells.Item(24, 4).Value = 1.031562099811399
$ws.Cells.Item(24, 5).Value = 1.024349743990518
$ws.Cells.Item(24, 6).Value = 1.038175914592582
$ws.Cells.Item(24, 9).Value = 1.029646278992668
$ws.Cells.Item(24, 10).Value = 1.029973835152332
$ws.Cells.Item(24, 11).Value = 1.034947721434242
$ws.Cells.Item(24, 12).Value = 1.027760970830239
$ws.Cells.Item(24, 13).Value = 1.041538400335805
$ws.Cells.Item(24, 14).Value = 1.014008753610247

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.025201405476731
$ws.Cells.Item(25, 4).Value = 1.032934577305623
$ws.Cells.Item(25, 5).Value = 1.025617203054928
$ws.Cells.Item(25, 6).Value = 1.039817591327277
$ws.Cells.Item(25, 9).Value = 1.029781887672494
$ws.Cells.Item(25, 10).Value = 1.030858782281178
$ws.Cells.Item(25, 11).Value = 1.035996139111766
$ws.Cells.Item(25, 12).Value = 1.028702081649963
$ws.Cells.Item(25, 13).Value = 1.042857547762743
$ws.Cells.Item(25, 14).Value = 1.014304936094387

Write-Host "applied 380 kV case updates"